$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 475-492 with revised weekly price data ---
# Row 475
$ws.Range("D475").Value = 45075
$ws.Range("J475").Value = 200
$ws.Range("K475").Value = 1100
$ws.Range("L475").Value = 1200
$ws.Range("M475").Value = 1150
$ws.Range("P475").Value = 1150

# Row 476
$ws.Range("D476").Value = 45075
$ws.Range("I476").Value = "Segunda"
$ws.Range("J476").Value = 150
$ws.Range("K476").Value = 1000
$ws.Range("L476").Value = 1000
$ws.Range("M476").Value = 1000
$ws.Range("P476").Value = 1000

# Row 477
$ws.Range("D477").Value = 44218
$ws.Range("I477").Value = "Primera"
$ws.Range("J477").Value = 2700
$ws.Range("K477").Value = 700
$ws.Range("L477").Value = 800
$ws.Range("M477").Value = 748
$ws.Range("P477").Value = 748

# Row 478
$ws.Range("D478").Value = 44921
$ws.Range("J478").Value = 400
$ws.Range("K478").Value = 700
$ws.Range("L478").Value = 800
$ws.Range("M478").Value = 750
$ws.Range("P478").Value = 750

# Row 479
$ws.Range("D479").Value = 44921
$ws.Range("I479").Value = "Segunda"
$ws.Range("J479").Value = 300
$ws.Range("K479").Value = 600
$ws.Range("L479").Value = 600
$ws.Range("M479").Value = 600
$ws.Range("P479").Value = 600

# Row 480
$ws.Range("D480").Value = 45014
$ws.Range("I480").Value = "Primera"
$ws.Range("J480").Value = 300
$ws.Range("K480").Value = 1000
$ws.Range("L480").Value = 1200
$ws.Range("M480").Value = 1100
$ws.Range("P480").Value = 1100

# Row 481
$ws.Range("D481").Value = 44648
$ws.Range("J481").Value = 160
$ws.Range("K481").Value = 750
$ws.Range("L481").Value = 800
$ws.Range("M481").Value = 775
$ws.Range("P481").Value = 775

# Row 482
$ws.Range("D482").Value = 44648
$ws.Range("J482").Value = 80

# Row 483
$ws.Range("D483").Value = 44469
$ws.Range("J483").Value = 300
$ws.Range("L483").Value = 750
$ws.Range("M483").Value = 725
$ws.Range("P483").Value = 725

# Row 484
$ws.Range("D484").Value = 44651
$ws.Range("J484").Value = 100
$ws.Range("K484").Value = 650
$ws.Range("L484").Value = 650
$ws.Range("M484").Value = 650
$ws.Range("P484").Value = 650

# Row 485
$ws.Range("D485").Value = 44160
$ws.Range("J485").Value = 2800
$ws.Range("K485").Value = 700
$ws.Range("L485").Value = 800
$ws.Range("M485").Value = 750
$ws.Range("P485").Value = 750

# Row 486
$ws.Range("D486").Value = 44160
$ws.Range("J486").Value = 1300
$ws.Range("K486").Value = 500
$ws.Range("L486").Value = 500
$ws.Range("M486").Value = 500
$ws.Range("P486").Value = 500

# Row 487
$ws.Range("D487").Value = 44988
$ws.Range("J487").Value = 300
$ws.Range("K487").Value = 900
$ws.Range("L487").Value = 900
$ws.Range("M487").Value = 900
$ws.Range("P487").Value = 900

# Row 488
$ws.Range("D488").Value = 44988
$ws.Range("I488").Value = "Segunda"
$ws.Range("J488").Value = 300
$ws.Range("L488").Value = 700
$ws.Range("M488").Value = 700
$ws.Range("P488").Value = 700

# Row 489
$ws.Range("D489").Value = 44999
$ws.Range("J489").Value = 500
$ws.Range("K489").Value = 1000
$ws.Range("M489").Value = 1100
$ws.Range("P489").Value = 1100

# Row 490
$ws.Range("D490").Value = 44939
$ws.Range("I490").Value = "Primera"
$ws.Range("J490").Value = 200
$ws.Range("K490").Value = 700
$ws.Range("L490").Value = 750
$ws.Range("M490").Value = 725
$ws.Range("P490").Value = 725

# Row 491
$ws.Range("D491").Value = 45040
$ws.Range("J491").Value = 300
$ws.Range("K491").Value = 1200
$ws.Range("L491").Value = 1200
$ws.Range("M491").Value = 1200
$ws.Range("O491").Value = "Región del Maule"
$ws.Range("P491").Value = 1200

# Row 492
$ws.Range("D492").Value = 45040
$ws.Range("I492").Value = "Segunda"
$ws.Range("J492").Value = 300
$ws.Range("K492").Value = 1000
$ws.Range("L492").Value = 1000
$ws.Range("M492").Value = 1000
$ws.Range("P492").Value = 1000

# --- Append new rows 493-495 with additional weekly price records ---
# Row 493
$ws.Range("A493").Value = 7
$ws.Range("B493").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C493").Value = "Ñuble"
$ws.Range("D493").Value = 44662
$ws.Range("E493").Value = 16
$ws.Range("F493").Value = 100112023
$ws.Range("G493").Value = "Brócoli"
$ws.Range("H493").Value = "Sin especificar"
$ws.Range("I493").Value = "Primera"
$ws.Range("J493").Value = 200
$ws.Range("K493").Value = 800
$ws.Range("L493").Value = 850
$ws.Range("M493").Value = 825
$ws.Range("N493").Value = "$/unidad"
$ws.Range("O493").Value = "Provincia de Diguillín"
$ws.Range("P493").Value = 825
$ws.Range("Q493").Value = 1
$ws.Range("R493").Value = "Hortaliza"
$ws.Range("D493").NumberFormat = $ws.Range("D474").NumberFormat

# Row 494
$ws.Range("A494").Value = 7
$ws.Range("B494").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C494").Value = "Ñuble"
$ws.Range("D494").Value = 44356
$ws.Range("E494").Value = 16
$ws.Range("F494").Value = 100112023
$ws.Range("G494").Value = "Brócoli"
$ws.Range("H494").Value = "Sin especificar"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 160
$ws.Range("K494").Value = 550
$ws.Range("L494").Value = 600
$ws.Range("M494").Value = 575
$ws.Range("N494").Value = "$/unidad"
$ws.Range("O494").Value = "Región del Maule"
$ws.Range("P494").Value = 575
$ws.Range("Q494").Value = 1
$ws.Range("R494").Value = "Hortaliza"
$ws.Range("D494").NumberFormat = $ws.Range("D474").NumberFormat

# Row 495
$ws.Range("A495").Value = 7
$ws.Range("B495").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C495").Value = "Ñuble"
$ws.Range("D495").Value = 44473
$ws.Range("E495").Value = 16
$ws.Range("F495").Value = 100112023
$ws.Range("G495").Value = "Brócoli"
$ws.Range("H495").Value = "Sin especificar"
$ws.Range("I495").Value = "Primera"
$ws.Range("J495").Value = 300
$ws.Range("K495").Value = 600
$ws.Range("L495").Value = 650
$ws.Range("M495").Value = 625
$ws.Range("N495").Value = "$/unidad"
$ws.Range("O495").Value = "Provincia de Diguillín"
$ws.Range("P495").Value = 625
$ws.Range("Q495").Value = 1
$ws.Range("R495").Value = "Hortaliza"
$ws.Range("D495").NumberFormat = $ws.Range("D474").NumberFormat

